$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in D/E hold numeric-looking strings (prices/percentages) that must
# stay text (matching the source inline-string cells), e.g. "1.000" or
# "94.50" would otherwise be auto-coerced to the numbers 1 / 94.5 and lose
# their trailing zeros. Mark just the cells we touch as Text first.

$cells = @("D2","E2","D3","E3","E4","D5","E5","D6","E6","D8","E8","E9","D10","E10","D11","E11","D12","E12","D13","E13","D14","E14","D15","E15","D16","E16","D17","E17","D18","E18","D19","E19","D20","E20","D21","E21","B22","C22","D22","E22","B23","C23","D23","E23","D24","E24","D25","E25","D26","E26","D27","E27","D28","E28","D29","E29","D30","D31","E31","D32","E32","D33","E33","D34","E34","D35","E35","D36","E36","D37","E37","D38","E38","D39","E39","D40","E40","D41","E41","D42","E42","D43","E43","D44","E44","D45","E45","B46","C46","D46","E46","B47","C47","D47","E47","D48","E48","B49","C49","D49","E49","B50","C50","D50","E50","B51","C51","D51","E51")
foreach ($addr in $cells) { $ws.Range($addr).NumberFormat = "@" }

$ws.Range("D2").Value = "30.303.22"
$ws.Range("E2").Value = "  +1.04%  "
$ws.Range("D3").Value = "1.922.27"
$ws.Range("E3").Value = "  +0.70%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "0.8165"
$ws.Range("E5").Value = "  +2.94%  "
$ws.Range("D6").Value = "244.48"
$ws.Range("E6").Value = "  +1.18%  "
$ws.Range("D8").Value = "0.3266"
$ws.Range("E8").Value = "  +3.45%  "
$ws.Range("E9").Value = "  +3.99%  "
$ws.Range("D10").Value = "0.07271"
$ws.Range("E10").Value = "  +5.28%  "
$ws.Range("D11").Value = "0.7963"
$ws.Range("E11").Value = "  +7.44%  "
$ws.Range("D12").Value = "0.08116"
$ws.Range("E12").Value = "  +1.49%  "
$ws.Range("D13").Value = "1.913.45"
$ws.Range("E13").Value = "  +0.26%  "
$ws.Range("D14").Value = "5.425"
$ws.Range("E14").Value = "  +4.57%  "
$ws.Range("D15").Value = "94.50"
$ws.Range("E15").Value = "  +1.70%  "
$ws.Range("D16").Value = "30.300.46"
$ws.Range("E16").Value = "  +1.01%  "
$ws.Range("D17").Value = "14.29"
$ws.Range("E17").Value = "  +2.47%  "
$ws.Range("D18").Value = "6.084"
$ws.Range("E18").Value = "  +3.85%  "
$ws.Range("D19").Value = "250.47"
$ws.Range("E19").Value = "  +2.01%  "
$ws.Range("D20").Value = "0.000007875"
$ws.Range("E20").Value = "  +1.66%  "
$ws.Range("D21").Value = "2.180.38"
$ws.Range("E21").Value = "  +1.05%  "
$ws.Range("B22").Value = "Chainlink"
$ws.Range("C22").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D22").Value = "8.207"
$ws.Range("E22").Value = "  +20.21%  "
$ws.Range("B23").Value = "Dai"
$ws.Range("C23").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D23").Value = "1.001"
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("D24").Value = "1.000"
$ws.Range("E24").Value = "  -0.02%  "
$ws.Range("D25").Value = "0.1675"
$ws.Range("E25").Value = "  +20.35%  "
$ws.Range("D26").Value = "9.524"
$ws.Range("E26").Value = "  +3.12%  "
$ws.Range("D27").Value = "167.96"
$ws.Range("E27").Value = "  +0.02%  "
$ws.Range("D28").Value = "19.09"
$ws.Range("E28").Value = "  +0.97%  "
$ws.Range("D29").Value = "2.157"
$ws.Range("E29").Value = "  +6.27%  "
$ws.Range("D30").Value = "1.371"
$ws.Range("D31").Value = "1.554"
$ws.Range("E31").Value = "  +2.72%  "
$ws.Range("D32").Value = "4.361"
$ws.Range("E32").Value = "  +1.29%  "
$ws.Range("D33").Value = "0.05730"
$ws.Range("E33").Value = "  +3.55%  "
$ws.Range("D34").Value = "4.155"
$ws.Range("E34").Value = "  +1.77%  "
$ws.Range("D35").Value = "1.307"
$ws.Range("E35").Value = "  +3.96%  "
$ws.Range("D36").Value = "0.7484"
$ws.Range("E36").Value = "  +2.32%  "
$ws.Range("D37").Value = "0.9998"
$ws.Range("E37").Value = "  +0.05%  "
$ws.Range("D38").Value = "2.728"
$ws.Range("E38").Value = "  +0.26%  "
$ws.Range("D39").Value = "0.01964"
$ws.Range("E39").Value = "  +2.20%  "
$ws.Range("D40").Value = "2.821"
$ws.Range("E40").Value = "  +1.36%  "
$ws.Range("D41").Value = "0.4514"
$ws.Range("E41").Value = "  +2.38%  "
$ws.Range("D42").Value = "74.76"
$ws.Range("E42").Value = "  +3.49%  "
$ws.Range("D43").Value = "5.995"
$ws.Range("E43").Value = "  -1.92%  "
$ws.Range("D44").Value = "0.8560"
$ws.Range("E44").Value = "  +2.28%  "
$ws.Range("D45").Value = "1.934"
$ws.Range("E45").Value = "  +3.24%  "
$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").Value = "1.040.22"
$ws.Range("E46").Value = "  +5.28%  "
$ws.Range("B47").Value = "PaxDollar"
$ws.Range("C47").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D47").Value = "1.0000"
$ws.Range("E47").Value = "  -0.03%  "
$ws.Range("D48").Value = "103.25"
$ws.Range("E48").Value = "  +2.77%  "
$ws.Range("B49").Value = "SynthetixNetwork"
$ws.Range("C49").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$ws.Range("D49").Value = "3.112"
$ws.Range("E49").Value = "  +11.21%  "
$ws.Range("B50").Value = "Aptos"
$ws.Range("C50").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D50").Value = "7.669"
$ws.Range("E50").Value = "  +1.89%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "9.949"
$ws.Range("E51").Value = "  +2.03%  "
